$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = "Death Shout"
